$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Updated")

# Update the DATE column (E2:E17) from 01/11/2019 (43770) to 01/02/2020 (43862)
$ws.Range("E2:E17").Value = 43862

# Update the selected cell / active selection on the sheet
$ws.Range("E12").Select()
